# Update cryptocurrency Price (column D) and Volume(1h) (column E) figures
# for rows 2-51 to reflect the latest scraped data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D cells to Text format before writing so values such as
# "1.005" or "15.60" are preserved exactly (not auto-converted to numbers).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.035.54'
$ws.Range('E2').Value = '  -2.12%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.664.80'
$ws.Range('E3').Value = '  -1.62%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.32'
$ws.Range('E5').Value = '  -0.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5092'
$ws.Range('E6').Value = '  +0.45%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.005'
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2628'
$ws.Range('E8').Value = '  -0.41%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06392'
$ws.Range('E9').Value = '  +2.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.66'
$ws.Range('E10').Value = '  -1.96%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07418'
$ws.Range('E11').Value = '  +1.56%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.670.54'
$ws.Range('E12').Value = '  -1.28%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.495'
$ws.Range('E13').Value = '  +0.45%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5803'
$ws.Range('E14').Value = '  -0.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008541'
$ws.Range('E15').Value = '  +3.32%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.19'
$ws.Range('E16').Value = '  -1.84%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.083.37'
$ws.Range('E17').Value = '  -2.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.897'
$ws.Range('E18').Value = '  -2.78%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.005'
$ws.Range('E19').Value = '  +0.08%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.71'
$ws.Range('E20').Value = '  -1.11%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '188.67'
$ws.Range('E21').Value = '  +1.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.194'
$ws.Range('E22').Value = '  -0.60%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.006'
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '145.79'
$ws.Range('E24').Value = '  +0.59%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '7.594'
$ws.Range('E25').Value = '  +0.44%  '
$ws.Range('E26').Value = '  +4.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.60'
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06560'
$ws.Range('E28').Value = '  +14.92%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.309'
$ws.Range('E29').Value = '  +0.98%  '
$ws.Range('E30').Value = '  -1.06%  '
$ws.Range('E31').Value = '  +0.75%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.498'
$ws.Range('E32').Value = '  -0.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.623'
$ws.Range('E33').Value = '  -1.67%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.017'
$ws.Range('E34').Value = '  +0.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6050'
$ws.Range('E35').Value = '  +1.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.366'
$ws.Range('E36').Value = '  -0.24%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.683'
$ws.Range('E37').Value = '  +0.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.200'
$ws.Range('E38').Value = '  +4.78%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01610'
$ws.Range('E39').Value = '  +0.91%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.074.65'
$ws.Range('E40').Value = '  -0.06%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8602'
$ws.Range('E41').Value = '  -0.39%  '
$ws.Range('E42').Value = '  +0.74%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.47'
$ws.Range('E43').Value = '  +2.18%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.812.47'
$ws.Range('E44').Value = '  -1.91%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000112'
$ws.Range('E45').Value = '  +5.20%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.09'
$ws.Range('E46').Value = '  -0.95%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.009'
$ws.Range('E47').Value = '  +0.32%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.006'
$ws.Range('E48').Value = '  -0.73%  '
$ws.Range('E49').Value = '  -0.09%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4290'
$ws.Range('E50').Value = '  -0.55%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.943'
$ws.Range('E51').Value = '  +4.14%  '
